$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-22 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-23 Friday", 2)
$d.Content.Find.Execute("941×7=", $true, $false, $false, $false, $false, $true, 1, $false, "196×9=", 2)
$d.Content.Find.Execute("254×9=", $true, $false, $false, $false, $false, $true, 1, $false, "454×8=", 2)
$d.Content.Find.Execute("296×4=", $true, $false, $false, $false, $false, $true, 1, $false, "894×4=", 2)
$d.Content.Find.Execute("999×7=", $true, $false, $false, $false, $false, $true, 1, $false, "914×2=", 2)
$d.Content.Find.Execute("827×8=", $true, $false, $false, $false, $false, $true, 1, $false, "953×2=", 2)
$d.Content.Find.Execute("102×6=", $true, $false, $false, $false, $false, $true, 1, $false, "880×5=", 2)
$d.Content.Find.Execute("242×2=", $true, $false, $false, $false, $false, $true, 1, $false, "356×3=", 2)
$d.Content.Find.Execute("866×2=", $true, $false, $false, $false, $false, $true, 1, $false, "645×9=", 2)
$d.Content.Find.Execute("428×2=", $true, $false, $false, $false, $false, $true, 1, $false, "407×3=", 2)
$d.Content.Find.Execute("883×7=", $true, $false, $false, $false, $false, $true, 1, $false, "369×5=", 2)
$d.Content.Find.Execute("955×8=", $true, $false, $false, $false, $false, $true, 1, $false, "441×9=", 2)
$d.Content.Find.Execute("482×2=", $true, $false, $false, $false, $false, $true, 1, $false, "482×9=", 2)
$d.Content.Find.Execute("468×8=", $true, $false, $false, $false, $false, $true, 1, $false, "975×4=", 2)
$d.Content.Find.Execute("315×6=", $true, $false, $false, $false, $false, $true, 1, $false, "611×8=", 2)
$d.Content.Find.Execute("165×2=", $true, $false, $false, $false, $false, $true, 1, $false, "297×4=", 2)
$d.Content.Find.Execute("769×5=", $true, $false, $false, $false, $false, $true, 1, $false, "569×5=", 2)
$d.Content.Find.Execute("622×4=", $true, $false, $false, $false, $false, $true, 1, $false, "942×2=", 2)
$d.Content.Find.Execute("524×3=", $true, $false, $false, $false, $false, $true, 1, $false, "866×4=", 2)
$d.Content.Find.Execute("482×5=", $true, $false, $false, $false, $false, $true, 1, $false, "125×3=", 2)
$d.Content.Find.Execute("889×3=", $true, $false, $false, $false, $false, $true, 1, $false, "226×8=", 2)
$d.Content.Find.Execute("634×2=", $true, $false, $false, $false, $false, $true, 1, $false, "341×7=", 2)
$d.Content.Find.Execute("825×8=", $true, $false, $false, $false, $false, $true, 1, $false, "421×5=", 2)
$d.Content.Find.Execute("525×9=", $true, $false, $false, $false, $false, $true, 1, $false, "734×5=", 2)
$d.Content.Find.Execute("666×9=", $true, $false, $false, $false, $false, $true, 1, $false, "478×9=", 2)
$d.Content.Find.Execute("370×9=", $true, $false, $false, $false, $false, $true, 1, $false, "445×4=", 2)
